$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '34.592.11'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +1.26%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.822.61'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.95%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.41%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '224.28'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.92%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.553'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.20%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.34%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '33.02'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +4.68%  '

# Row 9
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +3.35%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0706'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +6.84%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0931'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.04%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.066.86'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.06%  '

# Row 13
$ws.Range('B13').NumberFormat = "@"
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').NumberFormat = "@"
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.844.48'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +3.22%  '

# Row 14
$ws.Range('B14').NumberFormat = "@"
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').NumberFormat = "@"
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '11.20'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -1.01%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.651'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +3.26%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '34.613.21'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.35%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '4.32'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.97%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '69.62'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.34%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '252.99'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.62%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0799'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +7.74%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.25'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +7.15%  '

# Row 22
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.20%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.27'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.19%  '

# Row 24
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.17%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '161.78'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +3.51%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '16.51'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.08%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.18'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +2.14%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.115'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.34%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.996'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.47%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0528'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.90%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.82'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.41%  '

# Row 32
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.18%  '

# Row 33
$ws.Range('B33').NumberFormat = "@"
$ws.Range('B33').Value = 'Swop.fi'
$ws.Range('C33').NumberFormat = "@"
$ws.Range('C33').Value = 'https://coinranking.com/coin/yrCr2HW2c+swopfi-swop'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '503.45'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +866.05%  '

# Row 34
$ws.Range('B34').NumberFormat = "@"
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').NumberFormat = "@"
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.63'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +1.58%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.92'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +4.40%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.433.53'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -1.33%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.653'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +3.04%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.07'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.25%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0192'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +1.99%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.975'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +8.62%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '82.53'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -1.11%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.79'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -2.53%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.37'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.73%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.16'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +4.16%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '6.09'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +3.94%  '

# Row 46
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.59'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +5.63%  '

# Row 47
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.06'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.74%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0496'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -2.49%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.962.52'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.01%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '106.17'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +8.27%  '

# Row 51
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.10%  '
